$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "21/11/2024"
$ws.Range("B4").Value = "20/11/2024"
$ws.Range("B6").Value = "21/11/2024"
$ws.Range("B7").Value = "22/11/2024"
$ws.Range("B8").Value = "21/11/2024"
$ws.Range("B9").Value = "21/11/2024"
$ws.Range("B10").Value = "20/11/2024"
$ws.Range("B11").Value = "20/11/2024"
$ws.Range("B12").Value = "22/11/2024"
$ws.Range("B15").Value = "22/11/2024"
$ws.Range("B17").Value = "21/11/2024"
$ws.Range("B18").Value = "21/11/2024"
$ws.Range("B20").Value = "21/11/2024"
$ws.Range("B22").Value = "20/11/2024"
$ws.Range("B24").Value = "21/11/2024"
$ws.Range("B25").Value = "22/11/2024"
$ws.Range("B26").Value = "21/11/2024"
$ws.Range("B27").Value = "21/11/2024"
$ws.Range("B28").Value = "22/11/2024"
$ws.Range("B29").Value = "20/11/2024"
$ws.Range("B30").Value = "22/11/2024"
$ws.Range("B31").Value = "21/11/2024"
$ws.Range("B33").Value = "22/11/2024"
$ws.Range("B35").Value = "22/11/2024"
$ws.Range("B36").Value = "21/11/2024"
$ws.Range("B38").Value = "22/11/2024"
$ws.Range("B40").Value = "21/11/2024"
$ws.Range("B41").Value = "20/11/2024"
$ws.Range("B42").Value = "22/11/2024"
$ws.Range("B43").Value = "22/11/2024"
$ws.Range("B44").Value = "21/11/2024"
$ws.Range("B45").Value = "22/11/2024"
$ws.Range("B47").Value = "21/11/2024"
$ws.Range("B50").Value = "22/11/2024"
$ws.Range("B51").Value = "21/11/2024"
$ws.Range("B53").Value = "20/11/2024"
$ws.Range("B54").Value = "21/11/2024"
$ws.Range("B58").Value = "21/11/2024"
$ws.Range("B59").Value = "21/11/2024"
$ws.Range("B61").Value = "21/11/2024"
$ws.Range("B62").Value = "22/11/2024"
$ws.Range("B63").Value = "22/11/2024"
$ws.Range("B65").Value = "21/11/2024"
$ws.Range("B66").Value = "21/11/2024"
$ws.Range("B67").Value = "22/11/2024"
$ws.Range("B68").Value = "20/11/2024"
$ws.Range("B69").Value = "20/11/2024"
$ws.Range("B70").Value = "22/11/2024"
$ws.Range("B71").Value = "21/11/2024"
$ws.Range("B72").Value = "22/11/2024"
$ws.Range("B73").Value = "20/11/2024"
$ws.Range("B75").Value = "22/11/2024"
$ws.Range("B77").Value = "21/11/2024"
$ws.Range("B80").Value = "21/11/2024"
$ws.Range("B82").Value = "21/11/2024"
$ws.Range("B83").Value = "21/11/2024"
$ws.Range("B84").Value = "22/11/2024"
$ws.Range("B86").Value = "22/11/2024"
$ws.Range("B87").Value = "22/11/2024"
$ws.Range("B88").Value = "22/11/2024"
$ws.Range("B89").Value = "20/11/2024"
$ws.Range("B90").Value = "21/11/2024"
$ws.Range("B92").Value = "21/11/2024"
$ws.Range("B94").Value = "22/11/2024"
$ws.Range("B96").Value = "22/11/2024"
$ws.Range("B97").Value = "22/11/2024"
$ws.Range("B98").Value = "22/11/2024"
$ws.Range("B100").Value = "22/11/2024"
$ws.Range("B101").Value = "20/11/2024"
$ws.Range("B103").Value = "21/11/2024"
$ws.Range("B107").Value = "22/11/2024"
$ws.Range("B109").Value = "22/11/2024"
$ws.Range("B110").Value = "21/11/2024"
$ws.Range("B111").Value = "21/11/2024"
$ws.Range("B112").Value = "21/11/2024"
$ws.Range("B113").Value = "22/11/2024"
$ws.Range("B114").Value = "21/11/2024"
$ws.Range("B115").Value = "21/11/2024"
$ws.Range("B116").Value = "22/11/2024"
$ws.Range("B117").Value = "22/11/2024"
$ws.Range("B118").Value = "20/11/2024"
$ws.Range("B119").Value = "22/11/2024"
$ws.Range("B120").Value = "21/11/2024"
$ws.Range("B123").Value = "22/11/2024"
$ws.Range("B125").Value = "21/11/2024"
$ws.Range("B126").Value = "22/11/2024"
$ws.Range("B127").Value = "22/11/2024"
$ws.Range("B129").Value = "22/11/2024"
$ws.Range("B131").Value = "21/11/2024"
$ws.Range("B132").Value = "22/11/2024"
$ws.Range("B133").Value = "22/11/2024"
$ws.Range("B135").Value = "21/11/2024"
$ws.Range("B136").Value = "22/11/2024"
$ws.Range("B137").Value = "20/11/2024"
$ws.Range("B141").Value = "20/11/2024"
$ws.Range("B142").Value = "22/11/2024"
$ws.Range("B145").Value = "21/11/2024"
$ws.Range("B148").Value = "22/11/2024"
$ws.Range("B149").Value = "22/11/2024"
$ws.Range("B150").Value = "21/11/2024"
$ws.Range("B151").Value = "22/11/2024"
$ws.Range("B153").Value = "22/11/2024"
$ws.Range("B154").Value = "21/11/2024"
$ws.Range("B155").Value = "22/11/2024"
$ws.Range("B157").Value = "21/11/2024"
$ws.Range("B160").Value = "22/11/2024"
$ws.Range("B161").Value = "22/11/2024"
$ws.Range("B163").Value = "22/11/2024"
$ws.Range("B164").Value = "22/11/2024"
$ws.Range("B166").Value = "22/11/2024"
$ws.Range("B168").Value = "21/11/2024"
$ws.Range("B171").Value = "20/11/2024"
$ws.Range("B172").Value = "22/11/2024"
$ws.Range("B174").Value = "20/11/2024"
$ws.Range("B176").Value = "22/11/2024"
$ws.Range("B179").Value = "21/11/2024"
$ws.Range("B181").Value = "21/11/2024"
$ws.Range("B182").Value = "21/11/2024"
$ws.Range("B184").Value = "22/11/2024"
$ws.Range("B185").Value = "22/11/2024"
$ws.Range("B186").Value = "22/11/2024"
$ws.Range("B187").Value = "22/11/2024"
$ws.Range("B189").Value = "22/11/2024"
$ws.Range("B190").Value = "21/11/2024"
$ws.Range("B192").Value = "22/11/2024"
$ws.Range("B194").Value = "22/11/2024"
$ws.Range("B195").Value = "21/11/2024"
$ws.Range("B196").Value = "22/11/2024"
$ws.Range("B197").Value = "22/11/2024"
$ws.Range("B198").Value = "21/11/2024"
$ws.Range("B200").Value = "22/11/2024"
$ws.Range("B201").Value = "21/11/2024"
$ws.Range("B202").Value = "22/11/2024"
$ws.Range("B203").Value = "21/11/2024"
$ws.Range("B204").Value = "21/11/2024"
$ws.Range("B205").Value = "21/11/2024"
